# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
# Swap the data (columns B:AC) between pairs of rows whose underlying
# match records were reordered in the source feed. Column A (the
# display sequence number) stays where it is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(104, 105),
    @(115, 116),
    @(125, 126),
    @(181, 182),
    @(225, 226),
    @(238, 239)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
